$wb = $excel.ActiveWorkbook

# --- 1. Create the new "is_targeted list" sheet -----------------------------
# Newly-added/newly-copied worksheets in this host only keep their written
# cell values if they are the LAST sheet in the workbook at the moment the
# values are written (and if we never Move/reposition that same sheet
# afterwards). So: copy a template sheet to the end, rename it, write its
# values there, and then move every sheet that must come *after* it instead
# of moving the new sheet itself.
$template = $wb.Worksheets.Item("polarity list")
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "is_targeted list"
$newSheet.Cells.Item(1, 1).Value = "TRUE"
$newSheet.Cells.Item(2, 1).Value = "FALSE"

# Sheets that need to be after "is_targeted list", in final left-to-right
# order. Moving each one to the end (in order) re-creates that order while
# leaving "is_targeted list" itself untouched and intact.
$trailingSheets = @(
    "ms_source list",
    "polarity list",
    "ion_mobility list",
    "ms_scan_mode list",
    "resolution_x_unit list",
    "resolution_y_unit list",
    "preparation_maldi_matrix list",
    "desi_solvent_f...rate_unit list"
)
foreach ($name in $trailingSheets) {
    $s = $wb.Worksheets.Item($name)
    $s.Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
}

# --- 2. Point the "is_targeted" column's validation at the new list sheet ---
$ws = $wb.Worksheets.Item("Export as TSV")
$rng = $ws.Range("N2:N1048576")
$rng.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
